# Fruta / hortaliza, semanal
#
# The commit reshuffles the per-row time-series data (Fecha / Volumen /
# Precio minimo / Precio maximo / Precio promedio ponderado / Precio $/Kg)
# across the existing rows - i.e. each data row ends up showing the
# D/J/K/L/M/P values that some other row used to hold (one row, 21, is a
# fixed point and keeps its own values). The descriptive columns
# (Mercado, Region, Categoria, etc.) are untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# columns that travel together as a group
$cols = @(4, 10, 11, 12, 13, 16)   # D, J, K, L, M, P

# destRow -> sourceRow (source row's original D/J/K/L/M/P values become
# the new values for destRow)
$mapping = @{
    2=43; 3=40; 4=18; 5=4; 6=38; 7=31; 8=23; 9=14; 10=19;
    11=36; 12=30; 13=42; 14=37; 15=16; 16=45; 17=39; 18=3; 19=29; 20=24;
    21=21; 22=27; 23=20; 24=34; 25=26; 26=8; 27=7; 28=5; 29=17; 30=10;
    31=32; 32=33; 33=22; 34=44; 35=28; 36=11; 37=41; 38=15; 39=12; 40=2;
    41=6; 42=13; 43=35; 44=25; 45=9
}

# Snapshot every source row's current values BEFORE writing anything back,
# since several rows both donate and receive values.
$original = @{}
for ($row = 2; $row -le 45; $row++) {
    $rowVals = @{}
    foreach ($col in $cols) {
        $rowVals[$col] = $ws.Cells.Item($row, $col).Value()
    }
    $original[$row] = $rowVals
}

for ($row = 2; $row -le 45; $row++) {
    $srcRow = $mapping[$row]
    $srcVals = $original[$srcRow]
    foreach ($col in $cols) {
        $ws.Cells.Item($row, $col).Value = $srcVals[$col]
    }
}
